# Agency fiscal year downloaded data files - refresh Return A summary for
# LITTLESTOWN BOROUGH (previously GETTYSBURG BOROUGH) and update the print
# date stamp plus every offense count cell that changed with the new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report header: print date + agency name
$ws.Range("M3").Value = "Printed On: 10/23/2025"
$ws.Range("B16").Value = "Agency: LITTLESTOWN BOROUGH"

# Row 21: Rape(Total)
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1

# Row 22: Rape
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 0

# Row 23: Attempted Rape
$ws.Range("E23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1

# Row 24: Robbery(Total)
$ws.Range("E24").Value = 0
$ws.Range("I24").Value = 0

# Row 25: Robbery - Firearm
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 0

# Row 29: Assault(Total)
$ws.Range("E29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 0

# Row 30: Assault - Firearm
$ws.Range("E30").Value = 1
$ws.Range("I30").Value = 1

# Row 31: Assault - Knife or Cutting Instrument
$ws.Range("E31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0

# Row 32: Assault - Other Dangerous Weapon
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0

# Row 33: Assault - StrongArm(Hands,Fists,Feet,etc)
$ws.Range("E33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0

# Row 34: Burglary(Total)
$ws.Range("E34").Value = 2
$ws.Range("I34").Value = 2
$ws.Range("J34").Value = 1

# Row 36: Burglary - Unlawful Entry(No Force)
$ws.Range("E36").Value = 2
$ws.Range("I36").Value = 2
$ws.Range("J36").Value = 1

# Row 38: Larceny(Total) - Theft(Excluding Motor Vehicles)
$ws.Range("J38").Value = 22
$ws.Range("N38").Value = 2

# Row 43: Human Trafficking
$ws.Range("E43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0

# Row 45: Total Part I Offenses
$ws.Range("E45").Value = 33
$ws.Range("I45").Value = 33
$ws.Range("J45").Value = 24
$ws.Range("N45").Value = 2

# Row 46: Assault - Other(Simple, Not Aggravated)
$ws.Range("E46").Value = 9
$ws.Range("I46").Value = 9
$ws.Range("J46").Value = 12
$ws.Range("N46").Value = 2

# Row 47: Forgery and Counterfeiting
$ws.Range("E47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0

# Row 48: Fraud
$ws.Range("E48").Value = 5
$ws.Range("I48").Value = 5
$ws.Range("J48").Value = 5
$ws.Range("N48").Value = 0

# Row 49: Embezzlement
$ws.Range("E49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0

# Row 50: Stolen Property; Buying, Receiving, Possessing
$ws.Range("E50").Value = 2
$ws.Range("I50").Value = 2
$ws.Range("N50").Value = 1

# Row 51: Vandalism
$ws.Range("E51").Value = 6
$ws.Range("I51").Value = 6
$ws.Range("J51").Value = 6
$ws.Range("N51").Value = 1

# Row 52: Weapons; Carrying, Possessing, etc.
$ws.Range("E52").Value = 5
$ws.Range("I52").Value = 5
$ws.Range("J52").Value = 2

# Row 54: Sex Offenses (Except Line 2 and 16)
$ws.Range("E54").Value = 3
$ws.Range("I54").Value = 3
$ws.Range("J54").Value = 1
$ws.Range("N54").Value = 0

# Row 55: Drug Abuse Violations(Total)
$ws.Range("E55").Value = 24
$ws.Range("I55").Value = 24
$ws.Range("J55").Value = 26

# Row 61: Possession SubTotal
$ws.Range("E61").Value = 24
$ws.Range("I61").Value = 24
$ws.Range("J61").Value = 26

# Row 62: Opium - Cocaine
$ws.Range("E62").Value = 1
$ws.Range("I62").Value = 1

# Row 63: Marijuana
$ws.Range("E63").Value = 20
$ws.Range("I63").Value = 20
$ws.Range("J63").Value = 22

# Row 64: Synthetic
$ws.Range("E64").Value = 1
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = 1

# Row 65: Other
$ws.Range("E65").Value = 2
$ws.Range("I65").Value = 2
$ws.Range("J65").Value = 2

# Row 71: Driving Under The Influence
$ws.Range("E71").Value = 13
$ws.Range("I71").Value = 13
$ws.Range("J71").Value = 14

# Row 72: Liquor Laws
$ws.Range("E72").Value = 3
$ws.Range("I72").Value = 3
$ws.Range("J72").Value = 2
$ws.Range("N72").Value = 0

# Row 73: Drunkenness
$ws.Range("E73").Value = 2
$ws.Range("I73").Value = 2
$ws.Range("J73").Value = 1

# Row 74: Disorderly Conduct
$ws.Range("E74").Value = 4
$ws.Range("I74").Value = 4
$ws.Range("J74").Value = 4

# Row 76: All Other Offenses (Except Traffic)
$ws.Range("E76").Value = 14
$ws.Range("I76").Value = 14
$ws.Range("J76").Value = 13
$ws.Range("N76").Value = 0

# Row 77: Total Part II
$ws.Range("E77").Value = 90
$ws.Range("I77").Value = 90
$ws.Range("J77").Value = 88
$ws.Range("N77").Value = 6
